$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 ("Quản lí category sản phẩm" / "Admin, người bán"),
# shifting the rows below it up by one.
$ws.Rows("7:7").Delete()

# The "Quản lí loại sản phẩm" row (now row 9, previously row 10) also had
# its access target widened from "Admin" to "Admin, người bán".
$ws.Range("D9").Value = "Admin, người bán"

# Match the cursor position left behind after the edit.
$ws.Range("D9").Select() | Out-Null
